$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Columns("G").Delete() | Out-Null

$chartObjects = $ws.ChartObjects()
$chartObj1 = $chartObjects.Add(469.39, 335.25, 386.91, 216.0)
$chartObj1.Name = "Chart 1"
$chartObj1.Chart.ChartType = 51

$chartObj2 = $chartObjects.Add(902.05, 337.5, 418.96, 216.0)
$chartObj2.Name = "Chart 2"
$chartObj2.Chart.ChartType = 51
